$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New data rows (nper, pv, fv) for rows 2..32 - includes new RRI failure-case tests
$rowsData = @(
    @(0, 300, 400),
    @(0, -1, -3),
    @(1, -1, -3),
    @(12, 100, 10),
    @(12, 100, -90),
    @(5, 0, 0),
    @(5, -1, 5),
    @(5, 10, 10),
    @(2, 2, 8),
    @(2, 8, 2),
    @(2, 8, 0),
    @(2, 0, 10),
    @(12, -5, -6),
    @(1, -5, 0),
    @(12, -1, -1),
    @(12, 300, 300),
    @(12, 300, 400),
    @(12, 300, 4000),
    @(12, 300, 40000),
    @(24, 300, 400),
    @(24, 300, 4000),
    @(24, 300, 40000),
    @(38, 300, 400),
    @(38, 300, 4000),
    @(38, 300, 40000),
    @(8, 10000, 2441880),
    @(4, 5000, 6000),
    @(4, 5000, 10000),
    @(1, 250, 275),
    @(2, 250, 500),
    @(3, 250, 880)
)

$r = 2
foreach ($row in $rowsData) {
    $ws.Cells.Item($r, 1).Value = $row[0]
    $ws.Cells.Item($r, 2).Value = $row[1]
    $ws.Cells.Item($r, 3).Value = $row[2]
    $ws.Cells.Item($r, 4).Formula = "=_xlfn.RRI(A" + $r + ",B" + $r + ",C" + $r + ")"
    $r++
}

# Resize the Table1 ListObject to cover the new range
$lo = $ws.ListObjects.Item("Table1")
$lo.Resize($ws.Range("A1:D32"))

# Update the active selection
[void]$ws.Range("A5").Select()
